$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "custom accuracy": round the last data row's numeric values to 2 decimal places
$values = [ordered]@{
    "B5" = 8.65;  "C5" = 6.15;  "D5" = 0.85;  "E5" = 18.77; "F5" = 14.99;
    "G5" = 6.81;  "H5" = 29.23; "I5" = 10.47; "J5" = 4.56;  "K5" = 6.6;
    "L5" = 7.53;  "M5" = 7.89;  "N5" = 2.18;  "O5" = 6.77;  "P5" = 9.56;
    "Q5" = 5.88;  "R5" = 0.74;  "S5" = 0.49;  "T5" = 95.69; "U5" = 19.04;
    "V5" = 6.25;  "W5" = 12.61; "X5" = 6.55;  "Y5" = 1.22;  "Z5" = 13.78;
    "AA5" = 5.52; "AB5" = 4.99; "AC5" = 5.86; "AD5" = 7.87; "AE5" = 0.55;
    "AF5" = 26.77; "AG5" = 3.42; "AH5" = 7.81
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}

# "데이터 1000개" (1000 data points): drop the trailing extra sample row
$ws.Rows("6:6").Delete()
